# Update countries & provincias Spain
# Applies the data refresh described in the diff:
#  - Alemania (row 8): Casos activos / Recuperados updated
#  - India (row 19): Casos totales / Nuevos casos / Recuperados updated
#  - Rows 64/65 swap from Grecia/Hungria to Hungria/Grecia (sorted by Casos totales)
#    with refreshed daily figures
#  - Jamaica (row 118): Casos criticos updated

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Alemania
$ws.Range("D8").Value = 114500
$ws.Range("E8").Value = 37294

# India
$ws.Range("B19").Value = 27977
$ws.Range("C19").Value = 87
$ws.Range("E19").Value = 20570

# Row 64 becomes Hungria with refreshed numbers
$ws.Range("A64").Value = "Hungria"
$ws.Range("B64").Value = 2583
$ws.Range("C64").Value = 83
$ws.Range("D64").Value = 498
$ws.Range("E64").Value = 1805
$ws.Range("F64").Value = 61
$ws.Range("G64").Value = 8
$ws.Range("H64").Value = 280

# Row 65 becomes Grecia (previous row-64 figures)
$ws.Range("A65").Value = "Grecia"
$ws.Range("B65").Value = 2517
$ws.Range("C65").Value = 0
$ws.Range("D65").Value = 577
$ws.Range("E65").Value = 1806
$ws.Range("F65").Value = 46
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 134

# Jamaica
$ws.Range("F118").Value = 3
